$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: split the red run
#   "InnoDB的page刷到磁盘上要写4个操作系统block，在极端情况下(比如断电)不一定能保证4个块的写入原子性"
# into two runs, making the second half bold + underlined, and
# move the "_GoBack" bookmark so that it wraps exactly that second
# half.
# -----------------------------------------------------------------

$target = "在极端情况下(比如断电)不一定能保证4个块的写入原子性"

$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Replacement.ClearFormatting()
$rng1.Find.Replacement.Font.Bold = $true
$rng1.Find.Replacement.Font.Underline = 1
$found1 = $rng1.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)

if ($found1) {
    # $rng1 now spans exactly the replaced text; move the bookmark here.
    $d.Bookmarks.Add("_GoBack", $rng1)
}

# -----------------------------------------------------------------
# Edit 2: merge the two runs
#   "在完" + "成doublewrite页的写入后，再将doublewrite buffer中的页写入
#   各个表空间文件中。是否开启doublewrite还需要看具体情况。"
# back into a single run (this also removes the bookmark that used
# to sit between them, since it has already been relocated above).
# -----------------------------------------------------------------

$merged = "在完成doublewrite页的写入后，再将doublewrite buffer中的页写入各个表空间文件中。是否开启doublewrite还需要看具体情况。"

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$found2 = $rng2.Find.Execute($merged, $true, $false, $false, $false, $false, $true, 1, $false, $merged, 2)

Write-Host "edit1:" $found1 "edit2:" $found2
